# TC31_Canine_Filter_Breed-IrishSettr.xlsx
# "Fixed variables and query errors in Bread from TC30 to TC47"
#
# The CasesTab Cypher query (cell B2 on the "startup" sheet) referenced a
# `co` (cohort) variable/column that doesn't belong in this query - remove
# the trailing `OPTIONAL MATCH (co:cohort)...` RETURN column
# (`coalesce(co.cohort_description, '') AS `Cohort``) so the query only
# returns the columns that are actually produced by the MATCH/WITH clauses
# above it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

$fixedCasesQuery = @"
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)
WHERE demo.breed IN ['Irish Setter']
MATCH (c)<--(diag:diagnosis)
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (co:cohort)<-[*]-(c)
WITH DISTINCT c, s, demo, diag, co
RETURN  coalesce(c.case_id, '') AS ``Case ID`` ,
        coalesce(s.clinical_study_designation, '') AS ``Study Code`` ,
        coalesce(s.clinical_study_type, '') AS  ``Study Type``,
        coalesce(demo.breed, '') AS Breed ,
        coalesce(diag.disease_term, '') AS Diagnosis ,
        coalesce(diag.stage_of_disease, '') AS ``Stage of Disease`` ,
        coalesce(demo.patient_age_at_enrollment, '') AS Age ,
        coalesce(demo.sex, '') AS Sex ,
        coalesce(demo.neutered_indicator, '') AS ``Neutered Status``,
        coalesce(demo.weight, '') AS ``Weight (kg)``,
        coalesce(diag.best_response, '') AS ``Response to Treatment``
"@

$ws.Range("B2").Value = $fixedCasesQuery

# The author also left the selection sitting on B2 (scrolled back to the
# top of the sheet) instead of B4/row 4 when the file was resaved.
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
[void]$ws.Range("B2").Select()
